$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Rows 1-3: update existing single-value cells "100"/"0"/"74" -> "0M" ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# --- Insert 10 new single-cell rows right after row 3 (i.e. before the row
#     that is currently row 4), with values (top to bottom):
#     43, 0.00002, 0.00006, 0.00003, 0.00001, 0.00003, 0.00004, 0.00005,
#     0.00143, 100.0
#     Rows.Add(beforeRow) inserts immediately before "beforeRow", so to end
#     up with the values in forward (top-to-bottom) order we add them in
#     reverse, always re-fetching row 4 as the insertion point. ---
$newValues = @("43", "0.00002", "0.00006", "0.00003", "0.00001", "0.00003", "0.00004", "0.00005", "0.00143", "100.0")
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $beforeRow = $t.Rows.Item(4)
    $t.Rows.Add($beforeRow) | Out-Null
    $t.Rows.Item(4).Cells.Item(1).Range.Text = $newValues[$i]
}

# --- Collapse the final three rows (each a run of tab-separated values) down
#     to a single value. These rows used to be 34/35/36 (1-indexed) before the
#     10-row insertion above; they are now 44/45/46. ---
$rowCount = $t.Rows.Count
$t.Rows.Item($rowCount - 2).Cells.Item(1).Range.Text = "100"
$t.Rows.Item($rowCount - 1).Cells.Item(1).Range.Text = "0"
$t.Rows.Item($rowCount).Cells.Item(1).Range.Text = "74"
